$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "42.047.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.294.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "315.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "104.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.625"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.608"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "39.82"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0911"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "8.41"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.107"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.972"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "15.31"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "2.642.41"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "2.319.04"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "42.050.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "7.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "72.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "3.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "258.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "9.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "22.71"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "35.84"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "164.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.0887"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "5.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  +5.57%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.90"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.96%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0350"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "3.61"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "99.75"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +20.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "1.48"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "70.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "12.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "113.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "77.99"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.22%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "5.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("E51").Value = "  +2.55%  "
